# Applies the "Amelioration de la spredsheet wtf" edit:
#  - I6 label changes from "Temps pour faire 50$" to "Temps pour faire:"
#  - New cell J6 = 10, formatted as currency ("$"#,##0.00) - the value used to
#    be a hardcoded 10 inside the I7 formula
#  - I7 formula now references J6 instead of the hardcoded literal 10
#  - Selection moves to L13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in I6
$ws.Range("I6").Value = "Temps pour faire:"

# Update formula in I7 to reference J6 instead of the hardcoded 10
$ws.Range("I7").Formula = "=(J6/I5)/60"

# Add the new value cell J6 with a currency number format
$ws.Range("J6").Value = 10
$ws.Range("J6").NumberFormat = """$""#,##0.00"

# Update the active selection to L13, matching the saved sheet view state
$ws.Range("L13").Select()
